$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: delete row 6 (old "Juego de tronos. Realidades, ficciones, turismos" / AGBCEAAAQBAJ entry).
# This shifts rows 7-11 up by one (row 7 -> row 6, etc.).
$ws.Rows(6).Delete() | Out-Null

# Step 2: insert a new row at row 7 (right after the "Juego de tronos y la filosofia" row,
# which is now row 6) and populate it with the "Game of Thrones" pop-up book entry.
$ws.Rows(7).Insert() | Out-Null
$ws.Cells.Item(7, 1).Value = 'N12OngEACAAJ'
$ws.Cells.Item(7, 2).Value = 'Game of Thrones'
$ws.Cells.Item(7, 3).Value = 'Inspired by the Emmy® Award–winning credits sequence that opens each episode of the hit HBO® series, Game of Thrones: A Pop-Up Guide to Westeros is guaranteed to thrill the show’s legions of fans. Featuring stunning pop-up recreations of several key locations from the series, including the formidable castle of Winterfell, the lavish capital city King’s Landing, and the Wall’s stark majesty, this book—designed by renowned paper engineer Matthew Reinhart—takes you into the world of the series like never before. Game of Thrones: A Pop-Up Guide to Westeros features a total of five stunning spreads, which fold out to create a remarkable pop-up map of Westeros that is perfect for displaying. The book also contains numerous mini-pops that bring to life iconic elements of the show, such as direwolves, White Walkers, giants, and dragons. All the pops are accompanied by insightful text that relays the rich history of the Seven Kingdoms and beyond, forming a dynamic reference guide to the world of Game of Thrones. Visually spectacular and enthrallingly interactive, Game of Thrones: A Pop-Up Guide to Westeros sets a new standard for pop-up books and perfectly captures the epic scope and imagination of the series.'
# Column D holds a date-like string; force Text format first so Excel keeps it
# as the literal string "2014-06-10" instead of coercing it into a date serial.
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '2014-06-10'
$ws.Cells.Item(7, 5).Value = 'Matthew Reinhart'

# Step 3: delete the row holding the old "Juego de tronos" (Desconocido) stub entry
# (sbPj0AEACAAJ), which is now row 10 after the shifts above.
$ws.Rows(10).Delete() | Out-Null

# Step 4: insert a new row at row 11 (right after the "Realidades, ficciones, turismos"
# row, which is now row 10) and populate it with "The Complete Sherlock Holmes" entry.
$ws.Rows(11).Insert() | Out-Null
$ws.Cells.Item(11, 1).Value = 'zL5VdKN76gEC'
$ws.Cells.Item(11, 2).Value = 'The Complete Sherlock Holmes'
$ws.Cells.Item(11, 3).Value = 'Presents the four novels and fifty-six short stories which comprise the entire Sherlock Holmes saga'
# Column D holds a bare year number; force Text format first so Excel keeps it
# as the literal string "1930" instead of coercing it into a numeric value.
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '1930'
$ws.Cells.Item(11, 5).Value = 'Arthur Conan Doyle, Sir Arthur Conan Doyle'

